$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 984.5925999999999
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 984.5925999999999
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2953.7778
$ws.Range("N46").Value = -3191.7778
$ws.Range("M46").ClearContents()

# Row 60
$ws.Range("H60").Value = 984.5925999999999
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 984.5925999999999
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 2953.7778
$ws.Range("N60").Value = -3921.7778
$ws.Range("M60").ClearContents()

# Row 98
$ws.Range("H98").Value = 780.3461
$ws.Range("I98").Value = 828.875
$ws.Range("J98").Value = 198
$ws.Range("K98").Value = 828.875
$ws.Range("L98").Value = 198
$ws.Range("M98").Value = 669.125
$ws.Range("N98").Value = -3194

# Row 106
$ws.Range("H106").Value = 2603.8667
$ws.Range("I106").Value = 2668.9092
$ws.Range("J106").Value = 2425
$ws.Range("K106").Value = 2668.9092
$ws.Range("L106").Value = 2425
$ws.Range("M106").Value = -2037.9092
$ws.Range("N106").Value = -3687

# Row 122
$ws.Range("H122").Value = 780.3461
$ws.Range("I122").Value = 828.875
$ws.Range("J122").Value = 198
$ws.Range("K122").Value = 2486.625
$ws.Range("L122").Value = 594
$ws.Range("M122").Value = -36.625
$ws.Range("N122").Value = -5494

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 26103.73
$ws.Range("I32").Value = 5887.723
$ws.Range("J32").Value = 172108.22
$ws.Range("K32").Value = 5887.723
$ws.Range("L32").Value = 172108.22
$ws.Range("M32").Value = -5600.723
$ws.Range("N32").Value = -172682.22

# Row 61
$ws.Range("H61").Value = 1006.04254
$ws.Range("I61").Value = 920.23254
$ws.Range("J61").Value = 1928.5
$ws.Range("K61").Value = 920.23254
$ws.Range("L61").Value = 1928.5
$ws.Range("M61").Value = -708.23254
$ws.Range("N61").Value = -2352.5

# Row 74
$ws.Range("H74").Value = 720.2
$ws.Range("I74").Value = 740.43335
$ws.Range("J74").Value = 598.8
$ws.Range("K74").Value = 740.43335
$ws.Range("L74").Value = 598.8
$ws.Range("M74").Value = 133.56665
$ws.Range("N74").Value = -2346.8

# Row 77
$ws.Range("H77").Value = 720.2
$ws.Range("I77").Value = 740.43335
$ws.Range("J77").Value = 598.8
$ws.Range("K77").Value = 3702.16675
$ws.Range("L77").Value = 2994
$ws.Range("M77").Value = 665.8332499999997
$ws.Range("N77").Value = -11730

# Row 136
$ws.Range("H136").Value = 1006.04254
$ws.Range("I136").Value = 920.23254
$ws.Range("J136").Value = 1928.5
$ws.Range("K136").Value = 2760.69762
$ws.Range("L136").Value = 5785.5
$ws.Range("M136").Value = -210.6976199999999
$ws.Range("N136").Value = -10885.5

$ws = $wb.Worksheets.Item("BSM")
# Row 82
$ws.Range("H82").Value = 17311
$ws.Range("I82").Value = 3058.4
$ws.Range("J82").Value = 29188.166
$ws.Range("K82").Value = 3058.4
$ws.Range("L82").Value = 29188.166
$ws.Range("M82").Value = -2675.4
$ws.Range("N82").Value = -29954.166

# Row 85
$ws.Range("H85").Value = 17311
$ws.Range("I85").Value = 3058.4
$ws.Range("J85").Value = 29188.166
$ws.Range("K85").Value = 3058.4
$ws.Range("L85").Value = 29188.166
$ws.Range("M85").Value = -1732.4
$ws.Range("N85").Value = -31840.166

$ws = $wb.Worksheets.Item("CRP")
# Row 134
$ws.Range("H134").Value = 707.35297
$ws.Range("I134").Value = 557
$ws.Range("J134").Value = 1515.5
$ws.Range("K134").Value = 1671
$ws.Range("L134").Value = 4546.5
$ws.Range("M134").Value = 864
$ws.Range("N134").Value = -9616.5

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 6594.646
$ws.Range("I131").Value = 899.8570999999999
$ws.Range("J131").Value = 7042.551
$ws.Range("K131").Value = 2699.5713
$ws.Range("L131").Value = 21127.653
$ws.Range("M131").Value = 2340.4287
$ws.Range("N131").Value = -31207.653

# Row 133
$ws.Range("H133").Value = 4443.3
$ws.Range("I133").Value = 3600
$ws.Range("J133").Value = 5005.5
$ws.Range("K133").Value = 10800
$ws.Range("L133").Value = 15016.5
$ws.Range("M133").Value = -5740
$ws.Range("N133").Value = -25136.5

# Row 134
$ws.Range("H134").Value = 3260.476
$ws.Range("I134").Value = 2182.1428
$ws.Range("J134").Value = 3799.6428
$ws.Range("K134").Value = 6546.428400000001
$ws.Range("L134").Value = 11398.9284
$ws.Range("M134").Value = -1476.428400000001
$ws.Range("N134").Value = -21538.9284

# Row 136
$ws.Range("H136").Value = 933.3333
$ws.Range("I136").Value = 800
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 2400
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = 2700
$ws.Range("N136").Value = -16200

# Row 137
$ws.Range("H137").Value = 5057685
$ws.Range("I137").Value = 145372.86
$ws.Range("J137").Value = 7702776.5
$ws.Range("K137").Value = 436118.58
$ws.Range("L137").Value = 23108329.5
$ws.Range("M137").Value = -431018.58
$ws.Range("N137").Value = -23118529.5

# Row 138
$ws.Range("H138").Value = 10549.917
$ws.Range("I138").Value = 13066.556
$ws.Range("J138").Value = 3000
$ws.Range("K138").Value = 39199.66800000001
$ws.Range("L138").Value = 9000
$ws.Range("M138").Value = -34059.66800000001
$ws.Range("N138").Value = -19280

# Row 139
$ws.Range("H139").Value = 2554.48
$ws.Range("I139").Value = 1999
$ws.Range("J139").Value = 2924.8
$ws.Range("K139").Value = 5997
$ws.Range("L139").Value = 8774.400000000001
$ws.Range("M139").Value = -857
$ws.Range("N139").Value = -19054.4

# Row 140
$ws.Range("H140").Value = 5457.4614
$ws.Range("I140").Value = 7237.353
$ws.Range("J140").Value = 2095.4443
$ws.Range("K140").Value = 21712.059
$ws.Range("L140").Value = 6286.3329
$ws.Range("M140").Value = -16532.059
$ws.Range("N140").Value = -16646.3329

# Row 141
$ws.Range("H141").Value = 7939.353
$ws.Range("I141").Value = 9380.75
$ws.Range("J141").Value = 4480
$ws.Range("K141").Value = 28142.25
$ws.Range("L141").Value = 13440
$ws.Range("M141").Value = -22962.25
$ws.Range("N141").Value = -23800

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3269.077
$ws.Range("I7").Value = 1929.8
$ws.Range("J7").Value = 7733.3335
$ws.Range("K7").Value = 1929.8
$ws.Range("L7").Value = 7733.3335
$ws.Range("M7").Value = -1817.8
$ws.Range("N7").Value = -7957.3335

# Row 126
$ws.Range("H126").Value = 3269.077
$ws.Range("I126").Value = 1929.8
$ws.Range("J126").Value = 7733.3335
$ws.Range("K126").Value = 5789.4
$ws.Range("L126").Value = 23200.0005
$ws.Range("M126").Value = -3319.4
$ws.Range("N126").Value = -28140.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 574.7222
$ws.Range("I113").Value = 524.3077
$ws.Range("J113").Value = 603.2174
$ws.Range("K113").Value = 1572.9231
$ws.Range("L113").Value = 1809.6522
$ws.Range("M113").Value = 597.0769
$ws.Range("N113").Value = -6149.6522

# Row 136
$ws.Range("H136").Value = 449.92157
$ws.Range("I136").Value = 335.02222
$ws.Range("J136").Value = 1311.6666
$ws.Range("K136").Value = 1005.06666
$ws.Range("L136").Value = 3934.9998
$ws.Range("M136").Value = 1544.93334
$ws.Range("N136").Value = -9034.9998
